$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the error message cell (F4) and the last-update date cell (G4)
$ws.Range("F4").Value = "Database failed to get xai-xai CPI"
$ws.Range("G4").Value = "2022-09-07 22:25:55"
